# Apply updated Balance Sheet figures to the "DTE" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTE")

# Row 4 - Inventory
$ws.Range("B4").Value = 664000000.0
$ws.Range("C4").Value = 716000000.0
$ws.Range("D4").Value = 815000000.0
$ws.Range("E4").Value = 785000000.0
$ws.Range("F4").Value = 707000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 955000000.0
$ws.Range("C13").Value = 1029000000.0
$ws.Range("D13").Value = 964000000.0
$ws.Range("E13").Value = 1025000000.0
$ws.Range("F13").Value = 1028000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 3054000000.0
$ws.Range("C21").Value = 2984000000.0
$ws.Range("D21").Value = 2899000000.0
$ws.Range("E21").Value = 2813000000.0
$ws.Range("F21").Value = 2718000000.0

# Row 33 - Net Debt
$ws.Range("G33").Value = 17357000000.0

# Row 34 - Total Debt
$ws.Range("G34").Value = 17450000000.0
